$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" value (row 8, column B) ---
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2023-08-25T12:12:31-05:00"

# --- "Include from LOINC" sheet: replace concept row ---
# Before:
#   Row2: 11505-5 | Physician procedure note
#   Row3: 18842-5 | Discharge summary
# After:
#   Row2: 18842-5 | Discharge summary
#   Row3: 11506-3 | Progress note
$loincWs = $wb.Worksheets.Item("Include from LOINC")
$loincWs.Range("A2").Value = "18842-5"
$loincWs.Range("B2").Value = "Discharge summary"
$loincWs.Range("A3").Value = "11506-3"
$loincWs.Range("B3").Value = "Progress note"
